# Scheduled-runner update: refresh computed profit figures across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 33356.535
$ws.Range("I11").Value = 33356.535
$ws.Range("K11").Value = 33356.535
$ws.Range("M11").Value = -33216.535

$ws.Range("H32").Value = 5058.25
$ws.Range("I32").Value = 5058.4287
$ws.Range("J32").Value = 5058
$ws.Range("K32").Value = 5058.4287
$ws.Range("L32").Value = 5058
$ws.Range("M32").Value = -4732.4287
$ws.Range("N32").Value = -5710

$ws.Range("H43").Value = 2572.182
$ws.Range("J43").Value = 3026.3333
$ws.Range("L43").Value = 3026.3333
$ws.Range("N43").Value = -3164.3333

$ws.Range("H48").Value = 2259.5
$ws.Range("I48").Value = 2500
$ws.Range("J48").Value = 2019
$ws.Range("K48").Value = 7500
$ws.Range("L48").Value = 6057
$ws.Range("M48").Value = -7208
$ws.Range("N48").Value = -6641

$ws.Range("H56").Value = 2259.5
$ws.Range("I56").Value = 2500
$ws.Range("J56").Value = 2019
$ws.Range("K56").Value = 7500
$ws.Range("L56").Value = 6057
$ws.Range("M56").Value = -6966
$ws.Range("N56").Value = -7125

$ws.Range("H62").Value = 4195.8
$ws.Range("J62").Value = 4244.75
$ws.Range("L62").Value = 4244.75
$ws.Range("N62").Value = -5492.75

$ws.Range("H65").Value = 4195.8
$ws.Range("J65").Value = 4244.75
$ws.Range("L65").Value = 21223.75
$ws.Range("N65").Value = -27463.75

$ws.Range("I70").Value = 2100
$ws.Range("J70").Value = 1492
$ws.Range("K70").Value = 6300
$ws.Range("L70").Value = 4476
$ws.Range("M70").Value = -6030
$ws.Range("N70").Value = -5016

$ws.Range("I73").Value = 2100
$ws.Range("J73").Value = 1492
$ws.Range("K73").Value = 6300
$ws.Range("L73").Value = 4476
$ws.Range("M73").Value = -5364
$ws.Range("N73").Value = -6348

$ws.Range("H100").Value = 3812
$ws.Range("I100").Value = 2921.6667
$ws.Range("J100").Value = 5147.5
$ws.Range("K100").Value = 2921.6667
$ws.Range("L100").Value = 5147.5
$ws.Range("M100").Value = -2380.6667
$ws.Range("N100").Value = -6229.5

$ws.Range("H113").Value = 2941.4443
$ws.Range("I113").Value = 1996.2858
$ws.Range("K113").Value = 1996.2858
$ws.Range("M113").Value = 1257.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6999.423
$ws.Range("I2").Value = 7283.905
$ws.Range("J2").Value = 5804.6
$ws.Range("K2").Value = 7283.905
$ws.Range("L2").Value = 5804.6
$ws.Range("M2").Value = -7170.905
$ws.Range("N2").Value = -6030.6

$ws.Range("H45").Value = 1931.5
$ws.Range("I45").Value = 1717.5
$ws.Range("K45").Value = 1717.5
$ws.Range("M45").Value = -1340.5

$ws.Range("H102").Value = 1919.6
$ws.Range("I102").Value = 1781.909
$ws.Range("J102").Value = 2298.25
$ws.Range("K102").Value = 1781.909
$ws.Range("L102").Value = 2298.25
$ws.Range("M102").Value = -159.9090000000001
$ws.Range("N102").Value = -5542.25

$ws.Range("H110").Value = 8212.667
$ws.Range("I110").Value = 7917.9473
$ws.Range("K110").Value = 7917.9473
$ws.Range("M110").Value = -5872.9473

$ws.Range("H116").Value = 6999.423
$ws.Range("I116").Value = 7283.905
$ws.Range("J116").Value = 5804.6
$ws.Range("K116").Value = 7283.905
$ws.Range("L116").Value = 5804.6
$ws.Range("M116").Value = -4989.905
$ws.Range("N116").Value = -10392.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6999.423
$ws.Range("I3").Value = 7283.905
$ws.Range("J3").Value = 5804.6
$ws.Range("K3").Value = 7283.905
$ws.Range("L3").Value = 5804.6
$ws.Range("M3").Value = -7169.905
$ws.Range("N3").Value = -6032.6

$ws.Range("H107").Value = 15004.65
$ws.Range("I107").Value = 2482.353
$ws.Range("K107").Value = 2482.353
$ws.Range("M107").Value = -562.3530000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 726.6667
$ws.Range("I132").Value = 726.6667
$ws.Range("K132").Value = 2180.0001
$ws.Range("M132").Value = 349.9998999999998

$ws.Range("H134").Value = 3574275.8
$ws.Range("I134").Value = 2828.6365
$ws.Range("K134").Value = 8485.9095
$ws.Range("M134").Value = -5950.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 671.63635
$ws.Range("I44").Value = 417.4
$ws.Range("J44").Value = 883.5
$ws.Range("K44").Value = 1252.2
$ws.Range("L44").Value = 2650.5
$ws.Range("M44").Value = -854.1999999999998
$ws.Range("N44").Value = -3446.5

$ws.Range("H46").Value = 501000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H56").Value = 8336.733
$ws.Range("I56").Value = 8336.733
$ws.Range("K56").Value = 8336.733
$ws.Range("M56").Value = -7806.733

$ws.Range("H68").Value = 1148.1111
$ws.Range("I68").Value = 1249.6666
$ws.Range("J68").Value = 945
$ws.Range("K68").Value = 3748.9998
$ws.Range("L68").Value = 2835
$ws.Range("M68").Value = -2937.9998
$ws.Range("N68").Value = -4457

$ws.Range("H71").Value = 1148.1111
$ws.Range("I71").Value = 1249.6666
$ws.Range("J71").Value = 945
$ws.Range("K71").Value = 11246.9994
$ws.Range("L71").Value = 8505
$ws.Range("M71").Value = -7190.999400000001
$ws.Range("N71").Value = -16617

$ws.Range("H121").Value = 1501750.1
$ws.Range("I121").Value = 126311.375
$ws.Range("J121").Value = 3335668.5
$ws.Range("K121").Value = 378934.125
$ws.Range("L121").Value = 10007005.5
$ws.Range("M121").Value = -377624.125
$ws.Range("N121").Value = -10009625.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6033.3335
$ws.Range("I80").Value = 2625.5
$ws.Range("J80").Value = 7272.5454
$ws.Range("K80").Value = 2625.5
$ws.Range("L80").Value = 7272.5454
$ws.Range("M80").Value = -1627.5
$ws.Range("N80").Value = -9268.545399999999

$ws.Range("H83").Value = 6033.3335
$ws.Range("I83").Value = 2625.5
$ws.Range("J83").Value = 7272.5454
$ws.Range("K83").Value = 13127.5
$ws.Range("L83").Value = 36362.727
$ws.Range("M83").Value = -8135.5
$ws.Range("N83").Value = -46346.727

$ws.Range("H132").Value = 2481.3572
$ws.Range("I132").Value = 2477.0833
$ws.Range("K132").Value = 7431.249899999999
$ws.Range("M132").Value = -4901.249899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2703.7693
$ws.Range("I132").Value = 2348.0952
$ws.Range("K132").Value = 7044.285600000001
$ws.Range("M132").Value = -4514.285600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3647.25
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 3647.25
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H96").Value = 2685.889
$ws.Range("I96").Value = 4257.6665
$ws.Range("K96").Value = 4257.6665
$ws.Range("M96").Value = -2884.6665

$ws.Range("H126").Value = 3092.0557
$ws.Range("I126").Value = 2466.5715
$ws.Range("K126").Value = 7399.7145
$ws.Range("M126").Value = -4929.7145

